$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "NCTId" column (B) to host the
# new "status_label" field; this shifts NCTId..results from B:I to C:J.
$ws.Range("B1").EntireColumn.Insert()

# Header for the new column
$ws.Range("B1").Value = "status_label"

# String version of the "statut" emoji column (A), one value per data row
$ws.Range("B2").Value = "rouge"
$ws.Range("B3").Value = "vert"
$ws.Range("B4").Value = "orange"
$ws.Range("B5").Value = "rouge"
$ws.Range("B6").Value = "rouge"
$ws.Range("B7").Value = "rouge"
